# Updates cryptos list data: refreshed prices / 1h volume % and a few
# coins that swapped rank positions (row content updated in place).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the data range to Text so numeric-looking strings (e.g. "1.00",
# "0.999", thousand-dot prices) are kept verbatim instead of being
# parsed into numbers and losing formatting / precision.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "65.742.49"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "3.714.74"
$ws.Range("E3").Value = "  +4.66%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "411.06"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "133.86"
$ws.Range("E6").Value = "  +2.39%  "
$ws.Range("D7").Value = "3.708.12"
$ws.Range("E7").Value = "  +4.71%  "
$ws.Range("D8").Value = "0.629"
$ws.Range("E8").Value = "  -4.17%  "
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").Value = "0.738"
$ws.Range("E10").Value = "  -5.09%  "
$ws.Range("E11").Value = "  -4.51%  "
$ws.Range("D12").Value = "0.0000349"
$ws.Range("E12").Value = "  +13.16%  "
$ws.Range("D13").Value = "42.42"
$ws.Range("E13").Value = "  -1.68%  "
$ws.Range("D14").Value = "10.04"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "4.324.63"
$ws.Range("E15").Value = "  +5.13%  "
$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("D17").Value = "3.733.93"
$ws.Range("E17").Value = "  +5.28%  "
$ws.Range("D18").Value = "20.14"
$ws.Range("E18").Value = "  -1.80%  "
$ws.Range("D19").Value = "13.01"
$ws.Range("E19").Value = "  +4.13%  "
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("D21").Value = "66.206.07"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").Value = "422.63"
$ws.Range("E22").Value = "  -6.79%  "
$ws.Range("D23").Value = "14.93"
$ws.Range("E23").Value = "  +13.37%  "
$ws.Range("D24").Value = "87.10"
$ws.Range("E24").Value = "  -3.50%  "
$ws.Range("D25").Value = "3.03"
$ws.Range("E25").Value = "  -6.06%  "
$ws.Range("D26").Value = "36.38"
$ws.Range("E26").Value = "  +4.96%  "
$ws.Range("D27").Value = "3.20"
$ws.Range("E27").Value = "  -4.78%  "
$ws.Range("D28").Value = "9.53"
$ws.Range("E28").Value = "  -5.24%  "
$ws.Range("D29").Value = "5.18"
$ws.Range("E29").Value = "  +7.18%  "
$ws.Range("D30").Value = "12.52"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  +1.42%  "
$ws.Range("D32").Value = "2.73"
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("D33").Value = "7.03"
$ws.Range("E33").Value = "  -4.23%  "
$ws.Range("D34").Value = "41.47"
$ws.Range("E34").Value = "  +6.09%  "
$ws.Range("D35").Value = "0.160"
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("D36").Value = "55.88"
$ws.Range("E36").Value = "  -2.18%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "0.0471"
$ws.Range("E38").Value = "  -5.74%  "
$ws.Range("D39").Value = "2.98"
$ws.Range("E39").Value = "  +27.64%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0678"
$ws.Range("E40").Value = "  -15.57%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "0.141"
$ws.Range("E41").Value = "  -4.62%  "
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "3.39"
$ws.Range("E43").Value = "  +4.02%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "3.16"
$ws.Range("E44").Value = "  +23.22%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "26.81"
$ws.Range("E45").Value = "  +22.60%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "145.25"
$ws.Range("E46").Value = "  -2.57%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "2.09"
$ws.Range("E47").Value = "  +4.82%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "2.82"
$ws.Range("E48").Value = "  -6.52%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "4.24"
$ws.Range("E49").Value = "  -3.68%  "
$ws.Range("D50").Value = "2.53"
$ws.Range("E50").Value = "  -8.57%  "
$ws.Range("E51").Value = "  -5.65%  "

# Restore the original (default/general) cell formatting.
$dataRange.ClearFormats()
